$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style index 2) from A11 so the new date
# cells (A12:A15) get the same bold/bordered/centered date-number-format
# style used by the rest of column A, instead of minting a new style.
$ws.Range("A11").Copy()
$ws.Range("A12:A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 12 - Algyo / House
$ws.Range("A12").Value = 43845
$ws.Range("B12").Value = "Algyo"
$ws.Range("C12").Value = "House"
$ws.Range("D12").Value = 119.1052631578947
$ws.Range("E12").Value = 32676315.78947368
$ws.Range("F12").Value = 296798.3356605554
$ws.Range("G12").Value = 486.8421052631579
$ws.Range("H12").Value = 19

# Row 13 - Morahalom / House
$ws.Range("A13").Value = 43845
$ws.Range("B13").Value = "Morahalom"
$ws.Range("C13").Value = "House"
$ws.Range("D13").Value = 117.475
$ws.Range("E13").Value = 24720000
$ws.Range("F13").Value = 215462.9557050533
$ws.Range("G13").Value = 157.05
$ws.Range("H13").Value = 40

# Row 14 - Szeged / Garage (no land_size value)
$ws.Range("A14").Value = 43845
$ws.Range("B14").Value = "Szeged"
$ws.Range("C14").Value = "Garage"
$ws.Range("D14").Value = 17.72164948453608
$ws.Range("E14").Value = 4147010.30927835
$ws.Range("F14").Value = 241647.2934758712
$ws.Range("H14").Value = 97

# Row 15 - Szeged / House
$ws.Range("A15").Value = 43845
$ws.Range("B15").Value = "Szeged"
$ws.Range("C15").Value = "House"
$ws.Range("D15").Value = 161.6278195488722
$ws.Range("E15").Value = 55337731.82957394
$ws.Range("F15").Value = 989200.6957014774
$ws.Range("G15").Value = 466.5664160401003
$ws.Range("H15").Value = 798
